$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @('D2', '308.59'),
    @('E2', '0.10%'),
    @('D3', '41.10'),
    @('E3', '0.68%'),
    @('D4', '5.186'),
    @('E4', '1.25%'),
    @('D5', '0.07688'),
    @('E5', '0.80%'),
    @('D6', '1.646'),
    @('E6', '2.55%'),
    @('D7', '0.9146'),
    @('E7', '1.26%'),
    @('E8', '-1.62%'),
    @('D9', '0.1238'),
    @('E9', '10.42%'),
    @('E10', '2.44%'),
    @('D11', '0.09171'),
    @('E11', '-0.57%'),
    @('D12', '0.04223'),
    @('E12', '-1.29%'),
    @('D13', '0.1051'),
    @('E13', '-0.01%'),
    @('D14', '0.001255'),
    @('E14', '0.05%'),
    @('D15', '0.005747'),
    @('E15', '0.39%'),
    @('E16', '1,903.17%'),
    @('D17', '3.346'),
    @('E17', '-0.06%'),
    @('D18', '4.319'),
    @('E18', '1.61%'),
    @('D20', '7.414'),
    @('E20', '11.41%'),
    @('D21', '0.1383'),
    @('E21', '1.27%'),
    @('D22', '0.2713'),
    @('E22', '-3.18%'),
    @('D23', '0.04036'),
    @('E23', '-0.88%'),
    @('D24', '0.001268'),
    @('E24', '1.62%'),
    @('D25', '0.004071'),
    @('E25', '-1.11%'),
    @('D26', '0.0001301'),
    @('E26', '-0.04%'),
    @('D38', '0.02535'),
    @('E38', '6.44%'),
    @('D39', '0.05319'),
    @('E39', '2.64%'),
    @('D40', '0.007835'),
    @('E40', '0.59%'),
    @('E41', '1.03%'),
    @('D42', '0.006670'),
    @('E42', '-5.44%'),
    @('D43', '0.001861'),
    @('E43', '-4.65%'),
    @('D44', '0.008002'),
    @('E44', '0.74%'),
    @('D46', '0.00006707'),
    @('E46', '-4.29%'),
    @('E47', '-0.03%'),
    @('D48', '0.2939'),
    @('E48', '845.98%'),
    @('D50', '0.00002102'),
    @('E50', '-0.03%'),
    @('D51', '0.0002002'),
    @('E51', '-0.03%')
)

foreach ($u in $updates) {
    $ref = $u[0]
    $val = $u[1]
    $c = $ws.Range($ref)
    $s = $c.Style
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = $s
}

Write-Output "Applied $($updates.Count) cell updates"
